# edit.ps1
# Applies the "New crime data collected" update to the 46th Precinct CompStat
# weekly workbook: refreshes the report volume/number and date-range text,
# and updates the Murder / Rape rows (14-29) of the crime-complaints table
# with the newly collected weekly, 28-day, year-to-date and 2-year figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------
# xlPasteFormats = -4122
$xlPasteFormats = -4122

# Writes a *text* value into a cell (used for the "0" and "***.*" placeholder
# labels that appear instead of a number/percentage when a figure is not
# applicable). A leading apostrophe forces Excel to store the value as text
# even though it looks numeric; the cell format is then copied from a
# reference cell (row 30, which already holds a text-styled cell) so the
# resulting cell style matches the rest of the text cells in the table
# instead of picking up Excel's "quote prefix" style.
function Set-TextCell($cellAddr, $text) {
    $dst = $ws.Range($cellAddr)
    $dst.Value = "'" + $text
    $donor = $ws.Range("D30")
    $donor.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

# Writes a *numeric* value into a cell that previously held one of the text
# placeholders above, restoring the normal numeric cell style (copied from a
# reference cell of the right style in row 30: I30 for whole-number counts,
# K30 for percentage/decimal figures).
function Set-NumericCell($cellAddr, $num, $donorAddr) {
    $dst = $ws.Range($cellAddr)
    $dst.Value = $num
    $donor = $ws.Range($donorAddr)
    $donor.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# Report header: Volume/Number and the covered week's date range
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# ---------------------------------------------------------------------------
# Cells that change from the "n/a" text placeholder to a real number
# ---------------------------------------------------------------------------
    Set-NumericCell "D15" 2 "I30"
    Set-NumericCell "E15" 0 "K30"
    Set-NumericCell "D22" 1 "I30"
    Set-NumericCell "E22" -100 "K30"
    Set-NumericCell "D26" 2 "I30"
    Set-NumericCell "E26" 0 "K30"

# ---------------------------------------------------------------------------
# Cells that change from a real number to the "n/a" text placeholder
# ---------------------------------------------------------------------------
    Set-TextCell "C14" "0"
    Set-TextCell "G14" "0"
    Set-TextCell "H14" "***.*"
    Set-TextCell "C22" "0"
    Set-TextCell "C23" "0"
    Set-TextCell "C27" "0"
    Set-TextCell "C28" "0"
    Set-TextCell "D28" "0"
    Set-TextCell "E28" "***.*"
    Set-TextCell "C29" "0"
    Set-TextCell "D29" "0"
    Set-TextCell "E29" "***.*"

# ---------------------------------------------------------------------------
# Plain numeric value updates (cell style is unchanged)
# ---------------------------------------------------------------------------
    # Row 14
    $ws.Range("N14").Value = -75.806451612903

    # Row 15
    $ws.Range("C15").Value = 2
    $ws.Range("F15").Value = 4
    $ws.Range("H15").Value = 33.333333333333
    $ws.Range("I15").Value = 30
    $ws.Range("J15").Value = 25
    $ws.Range("K15").Value = 20
    $ws.Range("L15").Value = 50
    $ws.Range("M15").Value = 36.363636363636
    $ws.Range("N15").Value = -61.038961038961

    # Row 16
    $ws.Range("C16").Value = 5
    $ws.Range("D16").Value = 11
    $ws.Range("E16").Value = -54.545454545454
    $ws.Range("F16").Value = 26
    $ws.Range("G16").Value = 42
    $ws.Range("H16").Value = -38.095238095238
    $ws.Range("I16").Value = 369
    $ws.Range("J16").Value = 360
    $ws.Range("K16").Value = 2.5
    $ws.Range("L16").Value = 8.529411764705
    $ws.Range("M16").Value = -10.436893203883
    $ws.Range("N16").Value = -81.804733727810

    # Row 17
    $ws.Range("C17").Value = 15
    $ws.Range("D17").Value = 6
    $ws.Range("E17").Value = 150
    $ws.Range("F17").Value = 47
    $ws.Range("G17").Value = 37
    $ws.Range("H17").Value = 27.027027027027
    $ws.Range("I17").Value = 684
    $ws.Range("J17").Value = 571
    $ws.Range("K17").Value = 19.789842381786
    $ws.Range("L17").Value = 16.129032258064
    $ws.Range("M17").Value = 70.149253731343
    $ws.Range("N17").Value = -31.048387096774

    # Row 18
    $ws.Range("C18").Value = 4
    $ws.Range("D18").Value = 5
    $ws.Range("E18").Value = -20
    $ws.Range("F18").Value = 16
    $ws.Range("G18").Value = 24
    $ws.Range("H18").Value = -33.333333333333
    $ws.Range("I18").Value = 201
    $ws.Range("J18").Value = 175
    $ws.Range("K18").Value = 14.857142857142
    $ws.Range("L18").Value = -44.628099173553
    $ws.Range("M18").Value = 2.030456852791
    $ws.Range("N18").Value = -90.986547085201

    # Row 19
    $ws.Range("C19").Value = 15
    $ws.Range("E19").Value = 25
    $ws.Range("F19").Value = 39
    $ws.Range("G19").Value = 51
    $ws.Range("H19").Value = -23.529411764705
    $ws.Range("I19").Value = 566
    $ws.Range("J19").Value = 600
    $ws.Range("K19").Value = -5.666666666666
    $ws.Range("L19").Value = 3.096539162112
    $ws.Range("M19").Value = 96.527777777777
    $ws.Range("N19").Value = -38.276990185387

    # Row 20
    $ws.Range("C20").Value = 6
    $ws.Range("D20").Value = 4
    $ws.Range("E20").Value = 50
    $ws.Range("F20").Value = 10
    $ws.Range("H20").Value = -52.380952380952
    $ws.Range("I20").Value = 201
    $ws.Range("J20").Value = 158
    $ws.Range("K20").Value = 27.215189873417
    $ws.Range("L20").Value = 71.794871794871
    $ws.Range("M20").Value = 73.275862068965
    $ws.Range("N20").Value = -75.841346153846

    # Row 21
    $ws.Range("C21").Value = 47
    $ws.Range("D21").Value = 40
    $ws.Range("E21").Value = 17.5
    $ws.Range("F21").Value = 145
    $ws.Range("G21").Value = 178
    $ws.Range("H21").Value = -18.539325842696
    $ws.Range("I21").Value = 2066
    $ws.Range("J21").Value = 1901
    $ws.Range("K21").Value = 8.679642293529
    $ws.Range("L21").Value = 4.133064516129
    $ws.Range("M21").Value = 42.581090407177
    $ws.Range("N21").Value = -71.056318296441

    # Row 22
    $ws.Range("J22").Value = 36
    $ws.Range("K22").Value = -22.222222222222
    $ws.Range("L22").Value = -6.666666666666
    $ws.Range("M22").Value = -15.151515151515

    # Row 23
    $ws.Range("D23").Value = 3
    $ws.Range("E23").Value = -100
    $ws.Range("G23").Value = 5
    $ws.Range("H23").Value = -60
    $ws.Range("J23").Value = 44
    $ws.Range("K23").Value = -31.818181818181

    # Row 24
    $ws.Range("C24").Value = 21
    $ws.Range("D24").Value = 25
    $ws.Range("E24").Value = -16
    $ws.Range("F24").Value = 85
    $ws.Range("G24").Value = 90
    $ws.Range("H24").Value = -5.555555555555
    $ws.Range("I24").Value = 1263
    $ws.Range("J24").Value = 1128
    $ws.Range("K24").Value = 11.968085106383
    $ws.Range("L24").Value = 20.515267175572
    $ws.Range("M24").Value = 44.839449541284

    # Row 25
    $ws.Range("C25").Value = 18
    $ws.Range("D25").Value = 13
    $ws.Range("E25").Value = 38.461538461538
    $ws.Range("F25").Value = 61
    $ws.Range("G25").Value = 77
    $ws.Range("H25").Value = -20.779220779220
    $ws.Range("I25").Value = 1025
    $ws.Range("J25").Value = 831
    $ws.Range("K25").Value = 23.345367027677
    $ws.Range("L25").Value = 24.09200968523
    $ws.Range("M25").Value = -15.149006622516

    # Row 26
    $ws.Range("G26").Value = 3
    $ws.Range("H26").Value = 66.666666666666
    $ws.Range("I26").Value = 50
    $ws.Range("J26").Value = 40
    $ws.Range("K26").Value = 25
    $ws.Range("L26").Value = 38.888888888888

    # Row 27
    $ws.Range("D27").Value = 1
    $ws.Range("E27").Value = -100
    $ws.Range("F27").Value = 5
    $ws.Range("G27").Value = 7
    $ws.Range("H27").Value = -28.571428571428
    $ws.Range("J27").Value = 78
    $ws.Range("K27").Value = 7.692307692307
    $ws.Range("L27").Value = 31.25

    # Row 28
    $ws.Range("F28").Value = 2
    $ws.Range("G28").Value = 2
    $ws.Range("H28").Value = 0
    $ws.Range("L28").Value = -20.754716981132
    $ws.Range("N28").Value = -75.722543352601

    # Row 29
    $ws.Range("F29").Value = 2
    $ws.Range("G29").Value = 1
    $ws.Range("H29").Value = 100
    $ws.Range("L29").Value = -11.904761904761
    $ws.Range("N29").Value = -76.582278481012

